$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions): update F3 and F5
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2975
$ws1.Range("F5").Value = 73

# Sheet "全部类型" (All types, aggregate of all sheets): update F7 and F10
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 2975
$ws4.Range("F10").Value = 73
